$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "https://d5e0000019ce6eai"
$ws.Hyperlinks.Add($ws.Range("C1"), "https://d5e0000019ce6eai")
